$wb = $excel.ActiveWorkbook

# --- Update the "Status" text from "Ready for handoff" to "In Translation" ---
# This shared string is used by: Overview!E2, Overview!F2, zh-cn!C2, de-de!C2
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the "zh-cn"/"de-de" (Overview) and "Status" (per-locale) columns ---
# Target stored width (OOXML character units) is 13.4101845877511.
# The ColumnWidth COM property here rounds to the nearest 1/6 before Excel's
# standard 5/6-character padding is re-added when saved, so 12.5 is the input
# value that lands closest to the desired stored width.
$newColumnWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
